$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 73
$ws.Range("F2").Value = 50
$ws.Range("H2").Value = 62
$ws.Range("E4").Value = 29
$ws.Range("F4").Value = 10
$ws.Range("H4").Value = 24
$ws.Range("E5").Value = 10
$ws.Range("E7").Value = 18
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 7
$ws.Range("H12").Value = 7
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 194
$ws.Range("F15").Value = 113
$ws.Range("H15").Value = 154
$ws.Range("E17").Value = 149
$ws.Range("F17").Value = 80
$ws.Range("H17").Value = 112
$ws.Range("E18").Value = 146
$ws.Range("F18").Value = 73
$ws.Range("H18").Value = 110
$ws.Range("E19").Value = 78
$ws.Range("F19").Value = 48
$ws.Range("H19").Value = 61
$ws.Range("E24").Value = 30
$ws.Range("E25").Value = 31
$ws.Range("F25").Value = 20
$ws.Range("H25").Value = 28
$ws.Range("E27").Value = 21
$ws.Range("E28").Value = 26
$ws.Range("E34").Value = 32
$ws.Range("F34").Value = 13
$ws.Range("H34").Value = 16
$ws.Range("E37").Value = 70
$ws.Range("E38").Value = 95
$ws.Range("E40").Value = 30
$ws.Range("F40").Value = 21
$ws.Range("H40").Value = 23
$ws.Range("E44").Value = 36
$ws.Range("F44").Value = 22
$ws.Range("H44").Value = 32
$ws.Range("E48").Value = 48
$ws.Range("F48").Value = 33
$ws.Range("H48").Value = 39
$ws.Range("F49").Value = 48
$ws.Range("H49").Value = 65
$ws.Range("E51").Value = 15
$ws.Range("E52").Value = 11
$ws.Range("F52").Value = 7
$ws.Range("H52").Value = 7
$ws.Range("E57").Value = 23
$ws.Range("E59").Value = 12
$ws.Range("E61").Value = 38
$ws.Range("E62").Value = 61
$ws.Range("F62").Value = 20
$ws.Range("H62").Value = 34
$ws.Range("E63").Value = 52
$ws.Range("E64").Value = 44
$ws.Range("E65").Value = 43
$ws.Range("E70").Value = 55
$ws.Range("F70").Value = 30
$ws.Range("H70").Value = 43
$ws.Range("E71").Value = 51
$ws.Range("E72").Value = 56
$ws.Range("F72").Value = 31
$ws.Range("H72").Value = 42
$ws.Range("E73").Value = 37
$ws.Range("F73").Value = 19
$ws.Range("H73").Value = 31
$ws.Range("E74").Value = 22
$ws.Range("E75").Value = 22
$ws.Range("E79").Value = 50
$ws.Range("E82").Value = 20
$ws.Range("F82").Value = 9
$ws.Range("H82").Value = 15
$ws.Range("E89").Value = 51
$ws.Range("F89").Value = 25
$ws.Range("H89").Value = 31
